# Insert a new weekly price record as row 445 in the "Hortaliza, Feria
# Lagunitas de Puerto Montt - Cebollín" sheet, pushing the existing rows
# 445-529 down to 446-530 (dimension grows from A1:R529 to A1:R530).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 445..529 down by one to make room for the new record.
$ws.Rows.Item(445).Insert()

# Populate the newly inserted row 445 with the new record's data.
$ws.Cells.Item(445, 1).Value  = 4
$ws.Cells.Item(445, 2).Value  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(445, 3).Value  = 'Los Lagos'
$ws.Cells.Item(445, 4).Value  = 45211
$ws.Cells.Item(445, 5).Value  = 10
$ws.Cells.Item(445, 6).Value  = 100112037
$ws.Cells.Item(445, 7).Value  = 'Cebollín'
$ws.Cells.Item(445, 8).Value  = 'Sin especificar'
$ws.Cells.Item(445, 9).Value  = 'Primera'
$ws.Cells.Item(445, 10).Value = 70
$ws.Cells.Item(445, 11).Value = 6500
$ws.Cells.Item(445, 12).Value = 6500
$ws.Cells.Item(445, 13).Value = 6500
$ws.Cells.Item(445, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(445, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(445, 16).Value = 181
$ws.Cells.Item(445, 17).Value = 36
$ws.Cells.Item(445, 18).Value = 'Hortaliza'

# Copy the date cell's number format from the row above so the new row's
# date cell keeps the same "YYYY-MM-DD HH:MM:SS" style used throughout
# column D.
$ws.Cells.Item(444, 4).Copy()
$ws.Cells.Item(445, 4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(445, 4).Value = 45211
